# Add new transaction rows (4-9) to the worksheet, extending the table
# that currently spans A1:F3 up to A1:F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data block to append, one row per entry:
# Tipo, Descricao, Categoria, Valor, Vencimento (date), Status
$rows = @(
    @("Saida",   "Churrasco", "Lazer",        200,  "Pago"),
    @("Saida",   "Facul",     "Fixo",         400,  "Pago"),
    @("Saida",   "Mxrf11",    "Investimento", 90,   "Pago"),
    @("Entrada", "Mxrf11",    "Rendimentos",  0.9,  "Pago"),
    @("Entrada", "Salario",   "Pagamentos",   3000, "Pago"),
    @("Saida",   "Teste",     "Lazer",        50,   "Pendente")
)

# Same date used throughout column E of the existing rows: serial 46072
# (2026-02-19), kept with the same number format as E2/E3.
$dueDateSerial = 46072

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $eCell = $ws.Cells.Item($r, 5)
    # Reuse the exact same date format already applied to E2 (same numFmt
    # as the other Vencimento cells) instead of creating a new style.
    $eCell.NumberFormat = $ws.Cells.Item(2, 5).NumberFormat
    $eCell.Value = $dueDateSerial

    $ws.Cells.Item($r, 6).Value = $row[4]

    $r++
}
